# Applies the cell-value updates described by the commit diff to Sheet1
# of the active workbook. Only numeric values on rows 2, 5 and 7 change;
# no formatting/structure changes are needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "W2"  = 12
    "AC2" = 10
    "AD2" = 7.5
    "AE2" = 19
    "AS2" = 251
    "G5"  = 2.5
    "H5"  = 2.5
    "J5"  = 3.3
    "K5"  = 1.72
    "M5"  = 1.19
    "N5"  = 4.1
    "O5"  = 1.78
    "P5"  = 1.95
    "Q5"  = 3.25
    "R5"  = 1.3
    "S5"  = 1.75
    "T5"  = 2.02
    "U5"  = 2.42
    "V5"  = 1.5
    "W5"  = 5.1
    "X5"  = 10.25
    "Y5"  = 11
    "Z5"  = 29
    "AA5" = 32
    "AB5" = 65
    "AC5" = 4.1
    "AE5" = 22
    "AF5" = 200
    "AH5" = 6.5
    "AI5" = 16.5
    "AJ5" = 14
    "AM5" = 80
    "AN5" = 4.05
    "AO5" = 15
    "AP5" = 32
    "AQ5" = 80
    "AR5" = 175
    "AT5" = 2
    "AU5" = 8.5
    "AV5" = 120
    "AW5" = 5.1
    "AY5" = 40
    "G7"  = 2.3
    "I7"  = 2.7
    "Q7"  = 1.5
    "R7"  = 2.5
    "X7"  = 15
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
